# Upload new version with timestamp
# Fill in the pharmacy transactions table (rows 4-17), the totals row (18)
# and push the footer row down to row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data to insert -------------------------------------------------
# columns: A = م (seq no.), B = product name, H = ratio text, L = الرصيد الحالي, N = سعر البيع ratio
$data = @(
    @{A=1;  B="ANGIOFOX (EFFOX) 25MG LONG 30 CAPS.";        H="0:0";    L=114;   N=1},
    @{A=2;  B="AUGMENTIN 457MG/5ML SUSP. 70 ML";             H="1:0";    L=137;   N=1},
    @{A=3;  B="BLOKATENS 10/160MG 28 F.C.TABS.";             H="0:0";    L=160;   N=1},
    @{A=4;  B="COLOVATIL 30 F.C. TABS";                      H="0:0";    L=63;    N=1},
    @{A=5;  B="GAVISCON LIQUID 24 SACHETS 10 ML";            H="0:9";    L=12;    N=0.04},
    @{A=6;  B="GINKGO BILOBA 30 CAPS.";                      H="0:0";    L=186;   N=1},
    @{A=7;  B="MILGA ADVANCE 30 F.C. TABS";                  H="0:0";    L=136.5; N=1},
    @{A=8;  B="PERLOC 40MG 14 F.C.TAB.";                     H="0:0";    L=68.25; N=1},
    @{A=9;  B="RHINEX 0.05% INFANTILE NASAL DROPS 10 ML";    H="2:0";    L=18;    N=1},
    @{A=10; B="RIVO 320MG 20*10 TABS";                       H="1:2";    L=14.1;  N=0.1},
    @{A=11; B="VASTAREL MR 35MG 30 F.C.TAB.";                H="2:0";    L=175;   N=1},
    @{A=12; B="WATER FOR INJECTION AMP. 5 ML";               H="7816:0"; L=2.5;   N=1},
    @{A=13; B="سويت كوكو";                                   H="22:0";   L=25;    N=1},
    @{A=14; B="مرطب شفاه لونا جوز هند ابيض";                 H="3:0";    L=20;    N=1}
)

# row heights Excel settled on for each data row (alternating wrap heights)
$heights = @(24.75, 25.5, 24.75, 25.5, 25.5, 24.75, 25.5, 24.75, 25.5, 25.5, 24.75, 25.5, 24.75, 25.5)

# --- make room: row 4 already exists as the template row; insert 13 more
# rows below it so the table spans rows 4-17, pushing the totals row (old
# row 5) and the footer row (old row 6) down to rows 18 and 19.
$ws.Rows("5:17").Insert()

# --- clone row 4's formatting (font/fill/border/number format) onto the
# newly inserted rows 5-17, column group by column group (these groups
# are merged cells in row 4: B:G, H:K, L:M).
for ($r = 5; $r -le 17; $r++) {
    $ws.Range("A4").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("B4:G4").Copy()
    $ws.Range("B$r`:G$r").PasteSpecial(-4122)

    $ws.Range("H4:K4").Copy()
    $ws.Range("H$r`:K$r").PasteSpecial(-4122)

    $ws.Range("L4:M4").Copy()
    $ws.Range("L$r`:M$r").PasteSpecial(-4122)

    $ws.Range("N4").Copy()
    $ws.Range("N$r").PasteSpecial(-4122)

    $ws.Range("B$r`:G$r").Merge()
    $ws.Range("H$r`:K$r").Merge()
    $ws.Range("L$r`:M$r").Merge()
}
$excel.CutCopyMode = 0

# --- column B (product name) and column H (ratio text like "0:0") are
# stored as plain text -- force the Text number format before writing so
# values such as "0:0" / "1:2" are never re-interpreted.
$ws.Range("B4:B17").NumberFormat = "@"
$ws.Range("H4:H17").NumberFormat = "@"

# --- write the values row by row ------------------------------------
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 4 + $i
    $row = $data[$i]

    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("H$r").Value = $row.H
    $ws.Range("L$r").Value = $row.L
    $ws.Range("N$r").Value = $row.N

    $ws.Rows("$r`:$r").RowHeight = $heights[$i]
}

# --- totals row (old row 5, now row 18): set the grand total ---------
$ws.Range("K18").Value = 1131.35
$ws.Rows("18:18").RowHeight = 25.5

# --- footer row (old row 6, now row 19) keeps its original content; just
# restore its auto height ---------------------------------------------
$ws.Rows("19:19").RowHeight = 17.25

Write-Host "Inserted 14 transaction rows, updated totals and footer."
